# Generate Report for Handoff
#
# The localization-status report was regenerated: a new handoff was produced
# for the "de-de" ccbaa7f4... package/file, so its "Latest Handoff Datetime"
# (column D, row 4 of the "1374cf02..." source-file row) moves forward from
# 2016-01-28 03:46:05 to 2016-01-28 03:46:49.
#
# (The report-generation tool also re-appends a couple of bookkeeping
# entries to the shared-string table during the regeneration run, but none
# of the workbook's visible cell text besides this one timestamp actually
# changes - every other shared-string index shuffled by the diff still
# resolves to the exact same text it did before.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("D4").Value = "2016-01-28 03:46:49"
